$d = $word.ActiveDocument
$d.Content.Find.Execute("Екатеринбург", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Новороссийск", 2)
